# Update "想去人数" (column F) values on sheets "展览", "演出" and "全部类型"
# to reflect the regenerated data snapshot (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 282
$wsExhibit.Range("F4").Value = 43
$wsExhibit.Range("F5").Value = 3436
$wsExhibit.Range("F6").Value = 2154
$wsExhibit.Range("F7").Value = 416
$wsExhibit.Range("F8").Value = 166
$wsExhibit.Range("F9").Value = 52
$wsExhibit.Range("F10").Value = 39
$wsExhibit.Range("F11").Value = 1274
$wsExhibit.Range("F12").Value = 229
$wsExhibit.Range("F13").Value = 1614
$wsExhibit.Range("F14").Value = 115

# --- Sheet "演出" (performances) ---
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F3").Value = 1
$wsShow.Range("F4").Value = 2

# --- Sheet "全部类型" (all types) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 282
$wsAll.Range("F4").Value = 43
$wsAll.Range("F5").Value = 3436
$wsAll.Range("F6").Value = 2154
$wsAll.Range("F7").Value = 416
$wsAll.Range("F9").Value = 166
$wsAll.Range("F10").Value = 52
$wsAll.Range("F11").Value = 39
$wsAll.Range("F12").Value = 1
$wsAll.Range("F13").Value = 2
$wsAll.Range("F14").Value = 1274
$wsAll.Range("F15").Value = 229
$wsAll.Range("F16").Value = 1614
$wsAll.Range("F17").Value = 115
